$d = $word.ActiveDocument

# The page used to end with a couple of site-chrome paragraphs ("Ver no
# Jupiter Salvar em pdf Salvar em docx" and the "(c) 2020 ..." footer line)
# plus a spacer paragraph right after them. The rebuilt site no longer
# emits that chrome, so remove those paragraphs (and the now-redundant
# blank paragraph that trailed them), leaving the blank paragraph that
# already sat right after the final "... o aluno sera aprovado." text,
# followed directly by the page-break paragraph.

$jupiterPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Ver no Jupiter*") {
        $jupiterPara = $p
        break
    }
}

if ($jupiterPara -ne $null) {
    $copyrightPara = $jupiterPara.Next()
    $trailingBlankPara = $copyrightPara.Next()

    $startPos = $jupiterPara.Range.Start
    $endPos = $trailingBlankPara.Range.End

    $r = $d.Range($startPos, $endPos)
    $r.Delete()
}
